$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Script (Main)")

# Columns AF..AT for rows 3..6 get a new "ranking" style with thin left/right
# border and centered alignment, populated with 0 placeholders.
$cols = @("AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT")
foreach ($row in 3..6) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $cell.Value = 0
        $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
        $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
        $cell.HorizontalAlignment = -4108     # xlCenter
    }
}

$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("M59").Select()
